$wb = $excel.ActiveWorkbook

# Sheet "展览" and "全部类型" both hold the same event listing; update the
# "想去人数" (want-to-go count) figures for row 2 and row 5 on each.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1034
    $ws.Range("F5").Value = 440
}
